$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$caseIdQuery = @'
SELECT
    DISTINCT c.case_record_id AS "Case ID",
    st.clinical_study_designation AS "Study Code",
    st.clinical_study_type AS "Study Type",
    dmg.breed AS "Breed",
    diag.disease_term AS "Diagnosis",
    diag.stage_of_disease AS "Stage Of Disease",
    CASE 
    WHEN dmg.patient_age_at_enrollment = CAST(dmg.patient_age_at_enrollment AS INT) 
    THEN CAST(CAST(dmg.patient_age_at_enrollment AS INT) AS TEXT)
    ELSE CAST(dmg.patient_age_at_enrollment AS TEXT)
END AS "Age", 
    COALESCE(TRIM(dmg.sex), '') AS "Sex",   
    COALESCE(dmg.neutered_indicator, '') AS "Neutered Status",  
    COALESCE(
        CASE 
            WHEN dmg.weight = CAST(dmg.weight AS INT) THEN CAST(CAST(ROUND(dmg.weight, 0) AS INT) AS TEXT) 
            ELSE dmg.weight 
        END, 
    '') AS "Weight (kg)",  
    COALESCE(diag.best_response, '') AS "Response to Treatment",  
    COALESCE(coh.cohort_description, '') AS "Cohort"  
FROM 
    df_case c
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_sample smp ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_case_file cf ON cf."sample.sample_id" = smp.sample_id
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN
    df_cohort coh ON coh."study.clinical_study_designation" = st.clinical_study_designation
WHERE
    st.clinical_study_designation = 'COTC022' AND diag.best_response = 'Complete Response'
ORDER BY 
    c.case_record_id ASC
LIMIT 100;
'@

$studyFilesQuery = @'
SELECT DISTINCT
    sf.file_name AS "File Name",
    sf.file_type AS "File Type",
    'study' AS "Association",
    sf.file_description AS "Description",
    CASE
        WHEN sf.file_name LIKE '%.bai' THEN 'bai'
        WHEN sf.file_name LIKE '%.bam' THEN 'bam'
        WHEN sf.file_name LIKE '%.csv' THEN 'csv'
        WHEN sf.file_name LIKE '%.doc' THEN 'doc'
        WHEN sf.file_name LIKE '%.docx' THEN 'docx'
          WHEN cf.file_name LIKE '%.fastq.gz' THEN 'fastq'
        WHEN cf.file_name LIKE '%.vcf.gz' THEN 'vcf'
        WHEN cf.file_name LIKE '%.tsv.gz' THEN 'tsv'
        WHEN sf.file_name LIKE '%.pdf' THEN 'pdf'
        WHEN sf.file_name LIKE '%.rtf' THEN 'rtf'
        WHEN sf.file_name LIKE '%.tbi' THEN 'tbi'
        WHEN sf.file_name LIKE '%.tif' THEN 'tif'
        WHEN sf.file_name LIKE '%.xls' THEN 'xls'
        WHEN sf.file_name LIKE '%.xlsx' THEN 'xlsx'
        ELSE 'Unknown'
    END AS "Format",
      CASE     
    WHEN sf.file_size >= 1024 * 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 2) = CAST(ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'
            ELSE ROUND(sf.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
        END
    WHEN sf.file_size >= 1024 * 1024 THEN 
        CASE 
            WHEN ROUND(sf.file_size / (1024.0 * 1024.0), 2) = CAST(ROUND(sf.file_size / (1024.0 * 1024.0), 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'
            ELSE ROUND(sf.file_size / (1024.0 * 1024.0), 2) || ' MB'
        END
    WHEN sf.file_size >= 1024 THEN 
        CASE 
            WHEN ROUND(sf.file_size / 1024.0, 2) = CAST(ROUND(sf.file_size / 1024.0, 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'
            ELSE ROUND(sf.file_size / 1024.0, 2) || ' KB'
        END
    ELSE 
        CASE 
            WHEN ROUND(sf.file_size, 2) = CAST(ROUND(sf.file_size, 0) AS INT) 
            THEN CAST(CAST(ROUND(sf.file_size, 0) AS INT) AS TEXT) || ' Bytes'
            ELSE ROUND(sf.file_size, 2) || ' Bytes'
        END
END AS "Size",
    st.clinical_study_designation AS "Study Code"
FROM 
    df_case_file cf
JOIN 
    df_sample smp ON cf."sample.sample_id" = smp.sample_id
JOIN 
    df_case c ON smp."case.case_record_id" = c.case_record_id
JOIN 
    df_study st ON c."study.clinical_study_designation" = st.clinical_study_designation
JOIN 
    df_program p ON st."program.program_acronym" = p.program_acronym
JOIN 
    df_demographic dmg ON dmg."case.case_record_id" = c.case_record_id
JOIN 
    df_diagnosis diag ON diag."case.case_record_id" = c.case_record_id
JOIN 
    df_enrollment enr ON enr."case.case_record_id" = c.case_record_id
JOIN 
    df_publication pub ON pub."study.clinical_study_designation" = st.clinical_study_designation
LEFT JOIN 
    df_study_file sf ON sf."study.clinical_study_designation" = st.clinical_study_designation
WHERE
     st.clinical_study_designation = 'COTC022' AND diag.best_response = 'Complete Response'
ORDER BY 
    sf.file_name ASC
LIMIT 100;
'@

$ws.Cells.Item(2,2).Value = $caseIdQuery
$ws.Cells.Item(5,2).Value = $studyFilesQuery
$ws.Rows.Item(2).RowHeight = 409.6
$ws.Rows.Item(5).RowHeight = 409.6

$ws.Range("B5").Select()
